# "Generate Report for Archive"
#
# The localization status for this file moved from "Ready for handoff" to
# "In Translation". That status string shows up in three places:
#   - Overview sheet: columns "zh-cn" (E) and "de-de" (F), row 2
#   - zh-cn sheet: "Status" column (C), row 2
#   - de-de sheet: "Status" column (C), row 2
#
# Because the new text is shorter than the old text, the "Status"-ish
# columns that used to be autosized to fit "Ready for handoff" are
# narrowed to fit "In Translation" as well.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Narrow the columns that held the status text to fit the shorter string.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
